$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 34
$ws.Range("D2").Value = 108
$ws.Range("B3").Value = 28
$ws.Range("G3").Value = 15
$ws.Range("B4").Value = 15
$ws.Range("B5").Value = 21
$ws.Range("G5").Value = 9
$ws.Range("B6").Value = 55

$ws.Range("G6").Select()
